# Update cryptos list: price (D) and volume/1h (E) columns, plus two row
# identity swaps (B/C/D/E for rows 42, 43, 51) per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as literal text (e.g. "67.064.31"); Excel would
# otherwise auto-coerce plain-looking numbers ("599.42") to the Number
# type on assignment. Force Text format for the write, then clear the
# explicit format again so no stray style is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '67.054.18'
$ws.Range("E2").Value = '  +0.24%  '

Set-TextValue $ws.Range("D3") '3.493.19'
$ws.Range("E3").Value = '  -0.41%  '

$ws.Range("E4").Value = '  -0.03%  '

Set-TextValue $ws.Range("D5") '599.02'
$ws.Range("E5").Value = '  +0.50%  '

Set-TextValue $ws.Range("D6") '174.98'
$ws.Range("E6").Value = '  +3.14%  '

Set-TextValue $ws.Range("D7") '0.999'
$ws.Range("E7").Value = '  -0.06%  '

Set-TextValue $ws.Range("D8") '0.587'
$ws.Range("E8").Value = '  -0.58%  '

$ws.Range("E9").Value = '  -2.62%  '

Set-TextValue $ws.Range("D10") '7.14'
$ws.Range("E10").Value = '  -2.70%  '

$ws.Range("E11").Value = '  -1.33%  '

Set-TextValue $ws.Range("D12") '4.097.12'
$ws.Range("E12").Value = '  -0.47%  '

Set-TextValue $ws.Range("D13") '31.25'
$ws.Range("E13").Value = '  +10.06%  '

$ws.Range("E14").Value = '  +0.31%  '

Set-TextValue $ws.Range("D15") '67.054.85'
$ws.Range("E15").Value = '  +0.22%  '

$ws.Range("E16").Value = '  -2.10%  '

Set-TextValue $ws.Range("D17") '3.486.58'
$ws.Range("E17").Value = '  -0.12%  '

Set-TextValue $ws.Range("D18") '6.27'
$ws.Range("E18").Value = '  -1.22%  '

Set-TextValue $ws.Range("D19") '14.52'
$ws.Range("E19").Value = '  +2.62%  '

Set-TextValue $ws.Range("D20") '392.38'
$ws.Range("E20").Value = '  -0.92%  '

$ws.Range("E21").Value = '  -0.10%  '

Set-TextValue $ws.Range("D22") '73.29'
$ws.Range("E22").Value = '  -0.41%  '

$ws.Range("E23").Value = '  -0.08%  '

Set-TextValue $ws.Range("D24") '0.536'

Set-TextValue $ws.Range("D25") '5.70'
$ws.Range("E25").Value = '  +0.12%  '

$ws.Range("E26").Value = '  -2.27%  '

Set-TextValue $ws.Range("D27") '10.18'
$ws.Range("E27").Value = '  -0.48%  '

$ws.Range("E28").Value = '  -0.74%  '

$ws.Range("E29").Value = '  -0.36%  '

$ws.Range("E30").Value = '  -3.29%  '

$ws.Range("E31").Value = '  -3.43%  '

$ws.Range("E32").Value = '  -0.48%  '

$ws.Range("E33").Value = '  -1.83%  '

$ws.Range("E34").Value = '  -0.34%  '

$ws.Range("E35").Value = '  +1.04%  '

Set-TextValue $ws.Range("D36") '162.94'
$ws.Range("E36").Value = '  -0.75%  '

Set-TextValue $ws.Range("D37") '0.877'
$ws.Range("E37").Value = '  -2.19%  '

$ws.Range("E38").Value = '  +0.07%  '

Set-TextValue $ws.Range("D39") '7.04'
$ws.Range("E39").Value = '  +2.29%  '

$ws.Range("E40").Value = '  -1.77%  '

Set-TextValue $ws.Range("D41") '27.34'
$ws.Range("E41").Value = '  +1.01%  '

# Row 42: coin identity changed (re-sorted ranking)
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D42") '0.0730'
$ws.Range("E42").Value = '  -2.06%  '

# Row 43: coin identity changed (re-sorted ranking)
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D43") '26.06'
$ws.Range("E43").Value = '  -1.92%  '

Set-TextValue $ws.Range("D44") '2.799.92'
$ws.Range("E44").Value = '  -1.17%  '

$ws.Range("E45").Value = '  -0.91%  '

Set-TextValue $ws.Range("D46") '2.53'
$ws.Range("E46").Value = '  -4.07%  '

Set-TextValue $ws.Range("D47") '0.0300'
$ws.Range("E47").Value = '  -3.96%  '

Set-TextValue $ws.Range("D48") '337.73'
$ws.Range("E48").Value = '  -1.43%  '

Set-TextValue $ws.Range("D49") '1.07'
$ws.Range("E49").Value = '  -3.05%  '

Set-TextValue $ws.Range("D50") '33.81'
$ws.Range("E50").Value = '  -0.12%  '

# Row 51: coin identity changed (re-sorted ranking)
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range("D51") '0.846'
$ws.Range("E51").Value = '  -1.46%  '
